$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.382.56'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.774.39'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '351.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.24%  '
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +5.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.24'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.58%  '
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.82'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.214.55'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.795.82'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.923'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.436.97'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.47%  '
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.27'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.51%  '
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.45'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.44'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.74'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.87'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.26'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.90'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +8.78%  '
$ws.Range("E31").Value = '  -2.28%  '
$ws.Range("E32").Value = '  +7.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '52.02'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0444'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.54'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.21%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0837'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.54'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.09'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.89%  '
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.48'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.46'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.02'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.138.71'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("E48").Value = '  +5.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.223'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +17.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.47'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.894'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.55%  '
